$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin symbol list values (price, volume, hour) as scraped on 2023-02-07.
# Columns D, E, G on this sheet are stored as text, so force text number format
# before assigning values to avoid Excel auto-converting them to numbers/percentages.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D2").Value = "326.94"
$ws.Range("E2").Value = "0.10%"
$ws.Range("G2").Value = "5"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("B3").Value = "HuobiToken"
$ws.Range("C3").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D3").Value = "5.499"
$ws.Range("E3").Value = "-0.75%"
$ws.Range("G3").Value = "5"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("B4").Value = "Cronos"
$ws.Range("C4").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D4").Value = "0.08011"
$ws.Range("E4").Value = "-0.73%"
$ws.Range("G4").Value = "5"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("B5").Value = "FTXToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D5").Value = "1.989"
$ws.Range("E5").Value = "4.61%"
$ws.Range("G5").Value = "5"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.307"
$ws.Range("E6").Value = "-0.88%"
$ws.Range("G6").Value = "5"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "2.570"
$ws.Range("E7").Value = "-4.97%"
$ws.Range("G7").Value = "5"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9526"
$ws.Range("E8").Value = "0.56%"
$ws.Range("G8").Value = "5"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1127"
$ws.Range("E9").Value = "-4.40%"
$ws.Range("G9").Value = "5"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1871"
$ws.Range("E10").Value = "-1.27%"
$ws.Range("G10").Value = "5"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "10.61"
$ws.Range("E11").Value = "26.47%"
$ws.Range("G11").Value = "5"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09910"
$ws.Range("E12").Value = "-0.50%"
$ws.Range("G12").Value = "5"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.04579"
$ws.Range("E13").Value = "9.57%"
$ws.Range("G13").Value = "5"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.1066"
$ws.Range("E14").Value = "0.07%"
$ws.Range("G14").Value = "5"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001264"
$ws.Range("E15").Value = "-0.54%"
$ws.Range("G15").Value = "5"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04088"
$ws.Range("E16").Value = "-3.88%"
$ws.Range("G16").Value = "5"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005875"
$ws.Range("E17").Value = "-1.23%"
$ws.Range("G17").Value = "5"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("B18").Value = "OKB"
$ws.Range("C18").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D18").Value = "43.99"
$ws.Range("E18").Value = "-1.07%"
$ws.Range("G18").Value = "5"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("D19").Value = "3.357"
$ws.Range("E19").Value = "-6.72%"
$ws.Range("G19").Value = "5"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3477"
$ws.Range("E20").Value = "-0.28%"
$ws.Range("G20").Value = "5"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1410"
$ws.Range("E21").Value = "2.59%"
$ws.Range("G21").Value = "5"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2546"
$ws.Range("E22").Value = "-4.35%"
$ws.Range("G22").Value = "5"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001259"
$ws.Range("E23").Value = "1.59%"
$ws.Range("G23").Value = "5"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.57%"
$ws.Range("G24").Value = "5"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001159"
$ws.Range("E25").Value = "-5.96%"
$ws.Range("G25").Value = "5"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003744"
$ws.Range("E26").Value = "-6.51%"
$ws.Range("G26").Value = "5"

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "5"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "5"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "5"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "5"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "5"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "5"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "5"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "5"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "5"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "5"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "5"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02558"
$ws.Range("E38").Value = "-3.06%"
$ws.Range("G38").Value = "5"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05641"
$ws.Range("E39").Value = "1.76%"
$ws.Range("G39").Value = "5"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007521"
$ws.Range("E40").Value = "-2.15%"
$ws.Range("G40").Value = "5"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1396"
$ws.Range("E41").Value = "0.10%"
$ws.Range("G41").Value = "5"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007605"
$ws.Range("E42").Value = "-32.91%"
$ws.Range("G42").Value = "5"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002014"
$ws.Range("E43").Value = "-2.19%"
$ws.Range("G43").Value = "5"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008863"
$ws.Range("E44").Value = "1.88%"
$ws.Range("G44").Value = "5"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007112"
$ws.Range("E45").Value = "-0.07%"
$ws.Range("G45").Value = "5"

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.38%"
$ws.Range("G46").Value = "5"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("E47").Value = "54.86%"
$ws.Range("G47").Value = "5"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003091"
$ws.Range("E48").Value = "-9.79%"
$ws.Range("G48").Value = "5"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "-0.38%"
$ws.Range("G49").Value = "5"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "-0.38%"
$ws.Range("G50").Value = "5"

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "5"
